$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cosmed_xnat")

# New "Type" column (C) describing the data type of each XNAT field.
$ws.Range("C1").Value = "Type"
$ws.Range("C2").Value = "time"
$ws.Range("C3").Value = "float"
$ws.Range("C4").Value = "int"
$ws.Range("C5").Value = "int"
$ws.Range("C6").Value = "float"
$ws.Range("C7").Value = "float"
$ws.Range("C8").Value = "float"
$ws.Range("C9").Value = "float"
$ws.Range("C10").Value = "int"
$ws.Range("C11").Value = "int"
$ws.Range("C12").Value = "int"
$ws.Range("C13").Value = "float"
$ws.Range("C14").Value = "float"
$ws.Range("C15").Value = "float"
$ws.Range("C16").Value = "float"
$ws.Range("C17").Value = "float"
$ws.Range("C18").Value = "float"
$ws.Range("C19").Value = "int"
$ws.Range("C20").Value = "int"
$ws.Range("C21").Value = "int"

# Widen column B. This runtime snaps the stored column width to increments
# of 1/6 of a character, so 15.6 is the input that lands closest to the
# author's original (16.54296875) stored width. Also restore the selection
# like the author left it.
$ws.Columns.Item(2).ColumnWidth = 15.6
$ws.Range("C19:C21").Select()
